$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two rows of observations were dropped entirely from the data set:
#   "RM 232" (originally row 26) and "SC 92" (originally row 28).
# Deleting row 26 first shifts "SC 92" up to row 27, so we then delete row 27.
$ws.Rows(26).Delete()
$ws.Rows(27).Delete()

# After the two rows above are removed, the remaining rows shift up and
# occupy rows 2-33. Apply the updated / newly-missing values at their
# final (post-shift) row positions.
$ws.Range("D5").Value = ""
$ws.Range("E7").Value = ""
$ws.Range("D11").Value = -15.5
$ws.Range("F12").Value = ""
$ws.Range("F13").Value = 17.1
$ws.Range("F16").Value = 17.34
$ws.Range("F17").Value = 17.78
$ws.Range("C19").Value = 13.2
$ws.Range("D19").Value = ""
$ws.Range("F20").Value = 17.73
$ws.Range("C21").Value = ""
$ws.Range("F22").Value = 16.81
$ws.Range("C23").Value = 12.2
$ws.Range("D23").Value = -13.9
$ws.Range("E24").Value = -8.1
$ws.Range("F24").Value = ""
$ws.Range("D25").Value = -15.5
$ws.Range("F25").Value = ""
$ws.Range("B26").Value = ""
$ws.Range("B27").Value = -20.4
$ws.Range("C27").Value = ""
$ws.Range("D27").Value = ""
$ws.Range("E28").Value = -5.9
$ws.Range("F28").Value = ""
$ws.Range("B29").Value = ""
$ws.Range("D29").Value = ""
$ws.Range("F29").Value = ""
$ws.Range("E30").Value = ""
$ws.Range("F30").Value = ""
$ws.Range("E32").Value = ""
$ws.Range("C33").Value = 10.4
$ws.Range("D33").Value = -14.1
